# Update each "within 100" arithmetic answer cell in the single 20x5 table.
# Cells are addressed by (row, col) rather than by Find/Replace because one
# original equation text ("14+85=99") occurs twice but must become two
# different results depending on its position.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1; Col = 1; Old = '1+3=4'; New = '58-25=33' },
    @{ Row = 1; Col = 2; Old = '73+19=92'; New = '36-3=33' },
    @{ Row = 1; Col = 3; Old = '48-35=13'; New = '10+25=35' },
    @{ Row = 1; Col = 4; Old = '13+35=48'; New = '26-9=17' },
    @{ Row = 1; Col = 5; Old = '53+5=58'; New = '55-7=48' },
    @{ Row = 2; Col = 1; Old = '14+83=97'; New = '72-31=41' },
    @{ Row = 2; Col = 2; Old = '61+34=95'; New = '9+64=73' },
    @{ Row = 2; Col = 3; Old = '9+14=23'; New = '55+10=65' },
    @{ Row = 2; Col = 4; Old = '66-65=1'; New = '63+0=63' },
    @{ Row = 2; Col = 5; Old = '94+1=95'; New = '53+38=91' },
    @{ Row = 3; Col = 1; Old = '17+52=69'; New = '46+16=62' },
    @{ Row = 3; Col = 2; Old = '69+16=85'; New = '19+55=74' },
    @{ Row = 3; Col = 3; Old = '67-20=47'; New = '62-4=58' },
    @{ Row = 3; Col = 4; Old = '87-8=79'; New = '57+21=78' },
    @{ Row = 3; Col = 5; Old = '93-85=8'; New = '56+11=67' },
    @{ Row = 4; Col = 1; Old = '2+50=52'; New = '94-78=16' },
    @{ Row = 4; Col = 2; Old = '56-25=31'; New = '80+17=97' },
    @{ Row = 4; Col = 3; Old = '79-24=55'; New = '8+48=56' },
    @{ Row = 4; Col = 4; Old = '96-19=77'; New = '44+44=88' },
    @{ Row = 4; Col = 5; Old = '21+42=63'; New = '89-39=50' },
    @{ Row = 5; Col = 1; Old = '4+72=76'; New = '63+16=79' },
    @{ Row = 5; Col = 2; Old = '57-29=28'; New = '52+14=66' },
    @{ Row = 5; Col = 3; Old = '70-33=37'; New = '28-11=17' },
    @{ Row = 5; Col = 4; Old = '9+61=70'; New = '84-33=51' },
    @{ Row = 5; Col = 5; Old = '61-26=35'; New = '27+32=59' },
    @{ Row = 6; Col = 1; Old = '40-26=14'; New = '77+19=96' },
    @{ Row = 6; Col = 2; Old = '27-12=15'; New = '12+87=99' },
    @{ Row = 6; Col = 3; Old = '6+5=11'; New = '87-70=17' },
    @{ Row = 6; Col = 4; Old = '77-8=69'; New = '43-6=37' },
    @{ Row = 6; Col = 5; Old = '5+33=38'; New = '3+13=16' },
    @{ Row = 7; Col = 1; Old = '12+29=41'; New = '0+0=0' },
    @{ Row = 7; Col = 2; Old = '97-87=10'; New = '25+59=84' },
    @{ Row = 7; Col = 3; Old = '75+2=77'; New = '73-60=13' },
    @{ Row = 7; Col = 4; Old = '74-20=54'; New = '66-58=8' },
    @{ Row = 7; Col = 5; Old = '40+39=79'; New = '37+58=95' },
    @{ Row = 8; Col = 1; Old = '60-31=29'; New = '9+7=16' },
    @{ Row = 8; Col = 2; Old = '53+26=79'; New = '4-2=2' },
    @{ Row = 8; Col = 3; Old = '46+31=77'; New = '11+15=26' },
    @{ Row = 8; Col = 4; Old = '80+12=92'; New = '55-22=33' },
    @{ Row = 8; Col = 5; Old = '56-33=23'; New = '28+32=60' },
    @{ Row = 9; Col = 1; Old = '19+15=34'; New = '4+49=53' },
    @{ Row = 9; Col = 2; Old = '32+58=90'; New = '7+90=97' },
    @{ Row = 9; Col = 3; Old = '72+14=86'; New = '93-27=66' },
    @{ Row = 9; Col = 4; Old = '49-29=20'; New = '85-31=54' },
    @{ Row = 9; Col = 5; Old = '86-12=74'; New = '27+48=75' },
    @{ Row = 10; Col = 1; Old = '14+85=99'; New = '70-19=51' },
    @{ Row = 10; Col = 2; Old = '88-81=7'; New = '27+15=42' },
    @{ Row = 10; Col = 3; Old = '16+83=99'; New = '21+6=27' },
    @{ Row = 10; Col = 4; Old = '66+13=79'; New = '3+40=43' },
    @{ Row = 10; Col = 5; Old = '28-3=25'; New = '65-58=7' },
    @{ Row = 11; Col = 1; Old = '68-67=1'; New = '77-74=3' },
    @{ Row = 11; Col = 2; Old = '67-31=36'; New = '95-56=39' },
    @{ Row = 11; Col = 3; Old = '25-10=15'; New = '44-28=16' },
    @{ Row = 11; Col = 4; Old = '41+21=62'; New = '58+34=92' },
    @{ Row = 11; Col = 5; Old = '24+22=46'; New = '97-85=12' },
    @{ Row = 12; Col = 1; Old = '69+25=94'; New = '75-56=19' },
    @{ Row = 12; Col = 2; Old = '1+17=18'; New = '86-81=5' },
    @{ Row = 12; Col = 3; Old = '14+50=64'; New = '36-5=31' },
    @{ Row = 12; Col = 4; Old = '49+50=99'; New = '73-11=62' },
    @{ Row = 12; Col = 5; Old = '76-28=48'; New = '46+33=79' },
    @{ Row = 13; Col = 1; Old = '10+70=80'; New = '20+24=44' },
    @{ Row = 13; Col = 2; Old = '30-9=21'; New = '22-9=13' },
    @{ Row = 13; Col = 3; Old = '32+29=61'; New = '8+75=83' },
    @{ Row = 13; Col = 4; Old = '11+87=98'; New = '17+66=83' },
    @{ Row = 13; Col = 5; Old = '23+23=46'; New = '55+35=90' },
    @{ Row = 14; Col = 1; Old = '84-37=47'; New = '7+43=50' },
    @{ Row = 14; Col = 2; Old = '29+52=81'; New = '0+19=19' },
    @{ Row = 14; Col = 3; Old = '5+31=36'; New = '88-22=66' },
    @{ Row = 14; Col = 4; Old = '65-21=44'; New = '0+34=34' },
    @{ Row = 14; Col = 5; Old = '14+85=99'; New = '91+6=97' },
    @{ Row = 15; Col = 1; Old = '78+18=96'; New = '31-7=24' },
    @{ Row = 15; Col = 2; Old = '90+1=91'; New = '19+46=65' },
    @{ Row = 15; Col = 3; Old = '23+37=60'; New = '1+62=63' },
    @{ Row = 15; Col = 4; Old = '86-36=50'; New = '38-5=33' },
    @{ Row = 15; Col = 5; Old = '8+85=93'; New = '82-35=47' },
    @{ Row = 16; Col = 1; Old = '15+57=72'; New = '10+47=57' },
    @{ Row = 16; Col = 2; Old = '92-67=25'; New = '47-31=16' },
    @{ Row = 16; Col = 3; Old = '73-30=43'; New = '25+0=25' },
    @{ Row = 16; Col = 4; Old = '63-60=3'; New = '1+18=19' },
    @{ Row = 16; Col = 5; Old = '81-76=5'; New = '25+69=94' },
    @{ Row = 17; Col = 1; Old = '20+66=86'; New = '62-4=58' },
    @{ Row = 17; Col = 2; Old = '75-47=28'; New = '16+5=21' },
    @{ Row = 17; Col = 3; Old = '22+42=64'; New = '54-16=38' },
    @{ Row = 17; Col = 4; Old = '41+8=49'; New = '92-34=58' },
    @{ Row = 17; Col = 5; Old = '95-15=80'; New = '94-19=75' },
    @{ Row = 18; Col = 1; Old = '62-11=51'; New = '2+54=56' },
    @{ Row = 18; Col = 2; Old = '50+7=57'; New = '65+34=99' },
    @{ Row = 18; Col = 3; Old = '75-26=49'; New = '99-91=8' },
    @{ Row = 18; Col = 4; Old = '27+0=27'; New = '99-44=55' },
    @{ Row = 18; Col = 5; Old = '78-28=50'; New = '12+42=54' },
    @{ Row = 19; Col = 1; Old = '77-19=58'; New = '29+23=52' },
    @{ Row = 19; Col = 2; Old = '93-44=49'; New = '83-18=65' },
    @{ Row = 19; Col = 3; Old = '45-14=31'; New = '77-42=35' },
    @{ Row = 19; Col = 4; Old = '8+40=48'; New = '45-23=22' },
    @{ Row = 19; Col = 5; Old = '29+13=42'; New = '35-17=18' },
    @{ Row = 20; Col = 1; Old = '29-8=21'; New = '66-60=6' },
    @{ Row = 20; Col = 2; Old = '56-36=20'; New = '11+9=20' },
    @{ Row = 20; Col = 3; Old = '18+57=75'; New = '23-16=7' },
    @{ Row = 20; Col = 4; Old = '39-10=29'; New = '66+19=85' },
    @{ Row = 20; Col = 5; Old = '50-50=0'; New = '69-64=5' }
)

foreach ($rep in $replacements) {
    $cell = $t.Cell($rep.Row, $rep.Col)
    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $rep.Old) {
        Write-Host "MISMATCH at row $($rep.Row) col $($rep.Col): expected $($rep.Old) found $current"
    }
    $cell.Range.Text = $rep.New
}
